# Update odds/market cell values in Sheet1 to match the 2025-11-25 data refresh.
# Each block corresponds to one spreadsheet row (by A1-style row number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: F2, G2, Q2, W2
$ws.Cells.Item(2, 6).Value = 1.67
$ws.Cells.Item(2, 7).Value = 1.71
$ws.Cells.Item(2, 17).Value = 1.64
$ws.Cells.Item(2, 23).Value = 2.4

# Row 4: H4, I4, J4, K4, L4, T4, V4, AJ4, AL4, AN4
$ws.Cells.Item(4, 8).Value = 1.91
$ws.Cells.Item(4, 9).Value = 1.93
$ws.Cells.Item(4, 10).Value = 3.95
$ws.Cells.Item(4, 11).Value = 4
$ws.Cells.Item(4, 12).Value = 1.34
$ws.Cells.Item(4, 20).Value = 1.7
$ws.Cells.Item(4, 22).Value = 2.06
$ws.Cells.Item(4, 36).Value = 95
$ws.Cells.Item(4, 38).Value = 55
$ws.Cells.Item(4, 40).Value = 44

# Row 5: H5, I5, O5, P5, Q5, T5, V5, Y5, AA5, AI5, AL5, AO5
$ws.Cells.Item(5, 8).Value = 5.6
$ws.Cells.Item(5, 9).Value = 5.9
$ws.Cells.Item(5, 15).Value = 1.19
$ws.Cells.Item(5, 16).Value = 2.58
$ws.Cells.Item(5, 17).Value = 1.59
$ws.Cells.Item(5, 20).Value = 1.67
$ws.Cells.Item(5, 22).Value = 1.2
$ws.Cells.Item(5, 25).Value = 26
$ws.Cells.Item(5, 27).Value = 140
$ws.Cells.Item(5, 35).Value = 60
$ws.Cells.Item(5, 38).Value = 25
$ws.Cells.Item(5, 41).Value = 48

# Row 7: T7
$ws.Cells.Item(7, 20).Value = 1.55

# Row 8: F8, G8, K8, W8
$ws.Cells.Item(8, 6).Value = 2.14
$ws.Cells.Item(8, 7).Value = 2.34
$ws.Cells.Item(8, 11).Value = 3.3
$ws.Cells.Item(8, 23).Value = 1.75

# Row 9: H9
$ws.Cells.Item(9, 8).Value = 6.4

# Row 10: AM10
$ws.Cells.Item(10, 39).Value = 55

# Row 11: I11, K11, L11, N11, P11, Q11, S11
$ws.Cells.Item(11, 9).Value = 6.8
$ws.Cells.Item(11, 11).Value = 950
$ws.Cells.Item(11, 12).Value = 1.37
$ws.Cells.Item(11, 14).Value = 2.14
$ws.Cells.Item(11, 16).Value = 1.64
$ws.Cells.Item(11, 17).Value = 2.02
$ws.Cells.Item(11, 19).Value = 3.5

# Row 12: G12, H12, S12, W12
$ws.Cells.Item(12, 7).Value = 2.96
$ws.Cells.Item(12, 8).Value = 2.58
$ws.Cells.Item(12, 19).Value = 3.25
$ws.Cells.Item(12, 23).Value = 1.51

# Row 13: H13, I13, R13, AO13
$ws.Cells.Item(13, 8).Value = 2.6
$ws.Cells.Item(13, 9).Value = 2.64
$ws.Cells.Item(13, 18).Value = 1.44
$ws.Cells.Item(13, 41).Value = 19.5

# Row 14: F14, G14, J14, N14, Q14, R14, S14, T14, U14, W14, AG14, AN14
$ws.Cells.Item(14, 6).Value = 1.25
$ws.Cells.Item(14, 7).Value = 1.26
$ws.Cells.Item(14, 10).Value = 7.4
$ws.Cells.Item(14, 14).Value = 8.2
$ws.Cells.Item(14, 17).Value = 1.39
$ws.Cells.Item(14, 18).Value = 1.96
$ws.Cells.Item(14, 19).Value = 1.98
$ws.Cells.Item(14, 20).Value = 1.86
$ws.Cells.Item(14, 21).Value = 2.1
$ws.Cells.Item(14, 23).Value = 4.8
$ws.Cells.Item(14, 33).Value = 11.5
$ws.Cells.Item(14, 40).Value = 3.3

# Row 15: G15, H15, I15, N15, Q15, AA15
$ws.Cells.Item(15, 7).Value = 1.34
$ws.Cells.Item(15, 8).Value = 11.5
$ws.Cells.Item(15, 9).Value = 12
$ws.Cells.Item(15, 14).Value = 5.5
$ws.Cells.Item(15, 17).Value = 1.6
$ws.Cells.Item(15, 27).Value = 450

# Row 16: I16, J16, K16, P16, R16
$ws.Cells.Item(16, 9).Value = 4.7
$ws.Cells.Item(16, 10).Value = 3.85
$ws.Cells.Item(16, 11).Value = 3.95
$ws.Cells.Item(16, 16).Value = 2.24
$ws.Cells.Item(16, 18).Value = 1.49

# Row 17: F17, I17, J17, V17, Y17, Z17, AB17, AE17, AI17, AO17
$ws.Cells.Item(17, 6).Value = 3.2
$ws.Cells.Item(17, 9).Value = 2.3
$ws.Cells.Item(17, 10).Value = 3.9
$ws.Cells.Item(17, 22).Value = 1.77
$ws.Cells.Item(17, 25).Value = 14.5
$ws.Cells.Item(17, 26).Value = 17
$ws.Cells.Item(17, 28).Value = 17.5
$ws.Cells.Item(17, 31).Value = 21
$ws.Cells.Item(17, 35).Value = 29
$ws.Cells.Item(17, 41).Value = 12

# Row 18: G18, H18, I18, U18, AO18
$ws.Cells.Item(18, 7).Value = 2.48
$ws.Cells.Item(18, 8).Value = 2.84
$ws.Cells.Item(18, 9).Value = 2.94
$ws.Cells.Item(18, 21).Value = 3.15
$ws.Cells.Item(18, 41).Value = 12.5

# Row 19: N19
$ws.Cells.Item(19, 14).Value = 3.55
